$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values; regen save_data to use K instead of Strike#
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
